$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USNU")

# Insert a new column before column D (shifts D:K -> E:L)
$ws.Columns("D").Insert()

# Copy number formats from column E (which holds the old column D data) to the
# new column D so the new column matches the surrounding formatting.
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

Write-Output "Inserted column D and copied formats"
